$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New shared strings / header labels used by the new K:P result table
# ---------------------------------------------------------------------------
$corrDecalage = "Corr" + [char]0x00E9 + "lation et d" + [char]0x00E9 + "calage"
$corr          = "Corr" + [char]0x00E9 + "lation"
$manhDecalage  = "Manhattan et d" + [char]0x00E9 + "calage"

# ---------------------------------------------------------------------------
# 2) Column widths
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 10.59
$ws.Columns("K").ColumnWidth = 31.59
$ws.Columns("N").ColumnWidth = 20.88
$ws.Columns("P").ColumnWidth = 20.45

# ---------------------------------------------------------------------------
# 3) New header row (row 2), columns K:P mirroring A:I's headers
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = $ws.Range("A2").Value2
$ws.Range("L2").Value = $ws.Range("B2").Value2
$ws.Range("M2").Value = $corr
$ws.Range("N2").Value = $corrDecalage
$ws.Range("O2").Value = $ws.Range("D2").Value2
$ws.Range("P2").Value = $manhDecalage

# ---------------------------------------------------------------------------
# 4) New data rows, columns K:P
# ---------------------------------------------------------------------------
# Row 3
$ws.Range("K3").Value = $ws.Range("A3").Value2
$ws.Range("L3").Value = 0.48
$ws.Range("M3").Value = 0.52
$ws.Range("N3").Value = 0.48
$ws.Range("O3").Value = 0.52
$ws.Range("P3").Value = 0.44

# Row 4
$ws.Range("K4").Value = $ws.Range("A4").Value2
$ws.Range("L4").Value = 0.4204
$ws.Range("M4").Value = 0.63227276644352703
$ws.Range("N4").Value = 0.67878004829675198
$ws.Range("O4").Value = 0.473
$ws.Range("P4").Value = 0.4581

# Row 5
$ws.Range("K5").Value = $ws.Range("A5").Value2
$ws.Range("L5").Value = 0.80859999999999999
$ws.Range("M5").Value = 0.73043238586014902
$ws.Range("N5").Value = 0.79008896833391395
$ws.Range("O5").Value = 0.7024
$ws.Range("P5").Value = 0.6439

# Row 7
$ws.Range("K7").Value = $ws.Range("A7").Value2
$ws.Range("L7").Value = 0.72
$ws.Range("M7").Value = 0.76
$ws.Range("N7").Value = 0.92
$ws.Range("O7").Value = 0.73
$ws.Range("P7").Value = 0.96

# ---------------------------------------------------------------------------
# 5) Blank separator row (row 6) - touch every cell so it exists in sheetData
# ---------------------------------------------------------------------------
$ws.Range("A6:I6").Borders.LineStyle = 1
$ws.Range("K6:P6").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 6) Thin "all borders" box around every cell of the two tables
# ---------------------------------------------------------------------------
$ws.Range("A2:I5").Borders.LineStyle = 1
$ws.Range("A7:I7").Borders.LineStyle = 1
$ws.Range("K2:P5").Borders.LineStyle = 1
$ws.Range("K7:P7").Borders.LineStyle = 1

# H3:H5, H7 stay empty but still need the border (mirrors A:I template)
$ws.Range("H3:H5").Borders.LineStyle = 1
$ws.Range("H7").Borders.LineStyle = 1
